# Update cryptocurrency price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.400.14"
$ws.Range("D3").Value = "3.397.18"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.80"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.30"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.397.46"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.569"
$ws.Range("E9").Value = "  -8.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.24"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").Value = "  -4.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.422"
$ws.Range("E12").Value = "  -4.97%  "
$ws.Range("D13").Value = "3.978.71"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.80"
$ws.Range("E15").Value = "  -4.33%  "
$ws.Range("E16").Value = "  -9.82%  "
$ws.Range("D17").Value = "63.452.63"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "3.393.25"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.06"
$ws.Range("E19").Value = "  -5.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.46"
$ws.Range("E20").Value = "  -4.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.69"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.73"
$ws.Range("E22").Value = "  -4.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.82"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.512"
$ws.Range("E25").Value = "  -7.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000113"
$ws.Range("E26").Value = "  -5.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.64"
$ws.Range("E27").Value = "  -6.12%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.01"
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("E31").Value = "  -7.77%  "
$ws.Range("E32").Value = "  -2.85%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.75"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.89"
$ws.Range("E35").Value = "  -5.11%  "
$ws.Range("E36").Value = "  -7.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.06"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.835"
$ws.Range("E38").Value = "  +8.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.81"
$ws.Range("E39").Value = "  -5.16%  "
$ws.Range("D40").Value = "2.814.70"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.94"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0713"
$ws.Range("E43").Value = "  -6.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.35"
$ws.Range("E44").Value = "  -8.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.33"
$ws.Range("E45").Value = "  -6.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.46"
$ws.Range("E46").Value = "  -4.74%  "
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "325.25"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.30"
$ws.Range("E49").Value = "  +6.78%  "
$ws.Range("E50").Value = "  -5.36%  "
$ws.Range("E51").Value = "  -6.11%  "
